# Weekly fruit/vegetable price update: insert 3 new rows of data (week of
# 2022-01-24, serial 44585) for "Sandia" at "Vega Modelo de Temuco", just
# before the existing row 420, shifting the remaining rows down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 420-422; everything from the old row 420
# onward shifts down to 423 onward, preserving all of its data untouched.
$ws.Range("A420:A422").EntireRow.Insert()

# New row 420: Extra quality, Región del Maule
$ws.Range("A420").Value = 10
$ws.Range("B420").Value = "Vega Modelo de Temuco"
$ws.Range("C420").Value = "La Araucanía"
$ws.Range("D420").Value = 44585
$ws.Range("E420").Value = 9
$ws.Range("F420").Value = 100112028
$ws.Range("G420").Value = "Sandia"
$ws.Range("H420").Value = "Sin especificar"
$ws.Range("I420").Value = "Extra"
$ws.Range("J420").Value = 1500
$ws.Range("K420").Value = 3000
$ws.Range("L420").Value = 3000
$ws.Range("M420").Value = 3000
$ws.Range("N420").Value = "$/unidad"
$ws.Range("O420").Value = "Región del Maule"
$ws.Range("P420").Value = 3000
$ws.Range("Q420").Value = 1
$ws.Range("R420").Value = "Hortaliza"

# New row 421: Primera quality, Región del Maule
$ws.Range("A421").Value = 10
$ws.Range("B421").Value = "Vega Modelo de Temuco"
$ws.Range("C421").Value = "La Araucanía"
$ws.Range("D421").Value = 44585
$ws.Range("E421").Value = 9
$ws.Range("F421").Value = 100112028
$ws.Range("G421").Value = "Sandia"
$ws.Range("H421").Value = "Sin especificar"
$ws.Range("I421").Value = "Primera"
$ws.Range("J421").Value = 8000
$ws.Range("K421").Value = 2500
$ws.Range("L421").Value = 2500
$ws.Range("M421").Value = 2500
$ws.Range("N421").Value = "$/unidad"
$ws.Range("O421").Value = "Región del Maule"
$ws.Range("P421").Value = 2500
$ws.Range("Q421").Value = 1
$ws.Range("R421").Value = "Hortaliza"

# New row 422: Segunda quality, Región del Maule
$ws.Range("A422").Value = 10
$ws.Range("B422").Value = "Vega Modelo de Temuco"
$ws.Range("C422").Value = "La Araucanía"
$ws.Range("D422").Value = 44585
$ws.Range("E422").Value = 9
$ws.Range("F422").Value = 100112028
$ws.Range("G422").Value = "Sandia"
$ws.Range("H422").Value = "Sin especificar"
$ws.Range("I422").Value = "Segunda"
$ws.Range("J422").Value = 2000
$ws.Range("K422").Value = 2000
$ws.Range("L422").Value = 2000
$ws.Range("M422").Value = 2000
$ws.Range("N422").Value = "$/unidad"
$ws.Range("O422").Value = "Región del Maule"
$ws.Range("P422").Value = 2000
$ws.Range("Q422").Value = 1
$ws.Range("R422").Value = "Hortaliza"
